# Update NATMI LR-pair TPM-derived values for Psen1-Notch1 sheet
# Applies recomputed ligand/receptor/edge expression statistics (new TPM run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.539082
$ws.Range("H2").Value = 58.61724600000001
$ws.Range("I2").Value = 0.224220971665117
$ws.Range("J2").Value = 0.224220971665117
$ws.Range("M2").Value = 70.46291600000001
$ws.Range("N2").Value = 211.388748
$ws.Range("O2").Value = 0.5276750397950939
$ws.Range("P2").Value = 0.5276750397950939
$ws.Range("Q2").Value = 1376.780693683113
$ws.Range("R2").Value = 12391.02624314801
$ws.Range("S2").Value = 0.1183158101462852
$ws.Range("T2").Value = 0.1183158101462852
$ws.Range("G3").Value = 19.539082
$ws.Range("H3").Value = 58.61724600000001
$ws.Range("I3").Value = 0.224220971665117
$ws.Range("J3").Value = 0.224220971665117
$ws.Range("O3").Value = 0.07361176802536967
$ws.Range("P3").Value = 0.07361176802536967
$ws.Range("Q3").Value = 192.0637767603407
$ws.Range("R3").Value = 1728.573990843066
$ws.Range("S3").Value = 0.01650530215263557
$ws.Range("T3").Value = 0.01650530215263557
$ws.Range("G4").Value = 19.539082
$ws.Range("H4").Value = 58.61724600000001
$ws.Range("I4").Value = 0.224220971665117
$ws.Range("J4").Value = 0.224220971665117
$ws.Range("M4").Value = 42.505498
$ws.Range("N4").Value = 127.516494
$ws.Range("O4").Value = 0.3183105613832428
$ws.Range("P4").Value = 0.3183105613832428
$ws.Range("Q4").Value = 830.5184108728361
$ws.Range("R4").Value = 7474.665697855525
$ws.Range("S4").Value = 0.07137190336461956
$ws.Range("T4").Value = 0.07137190336461956
$ws.Range("G5").Value = 19.539082
$ws.Range("H5").Value = 58.61724600000001
$ws.Range("I5").Value = 0.224220971665117
$ws.Range("J5").Value = 0.224220971665117
$ws.Range("M5").Value = 10.73653933333333
$ws.Range("N5").Value = 32.209618
$ws.Range("O5").Value = 0.08040263079629371
$ws.Range("P5").Value = 0.08040263079629371
$ws.Range("Q5").Value = 209.7821224302254
$ws.Range("R5").Value = 1888.039101872028
$ws.Range("S5").Value = 0.01802795600157663
$ws.Range("T5").Value = 0.01802795600157663
$ws.Range("I6").Value = 0.3010605798326856
$ws.Range("J6").Value = 0.3010605798326856
$ws.Range("M6").Value = 70.46291600000001
$ws.Range("N6").Value = 211.388748
$ws.Range("O6").Value = 0.5276750397950939
$ws.Range("P6").Value = 0.5276750397950939
$ws.Range("Q6").Value = 1848.597795578859
$ws.Range("R6").Value = 16637.38016020973
$ws.Range("S6").Value = 0.1588621534439464
$ws.Range("T6").Value = 0.1588621534439464
$ws.Range("I7").Value = 0.3010605798326856
$ws.Range("J7").Value = 0.3010605798326856
$ws.Range("O7").Value = 0.07361176802536967
$ws.Range("P7").Value = 0.07361176802536967
$ws.Range("S7").Value = 0.02216160156422694
$ws.Range("T7").Value = 0.02216160156422694
$ws.Range("I8").Value = 0.3010605798326856
$ws.Range("J8").Value = 0.3010605798326856
$ws.Range("M8").Value = 42.505498
$ws.Range("N8").Value = 127.516494
$ws.Range("O8").Value = 0.3183105613832428
$ws.Range("P8").Value = 0.3183105613832428
$ws.Range("Q8").Value = 1115.133666945909
$ws.Range("R8").Value = 10036.20300251318
$ws.Range("S8").Value = 0.09583076217690675
$ws.Range("T8").Value = 0.09583076217690675
$ws.Range("I9").Value = 0.3010605798326856
$ws.Range("J9").Value = 0.3010605798326856
$ws.Range("M9").Value = 10.73653933333333
$ws.Range("N9").Value = 32.209618
$ws.Range("O9").Value = 0.08040263079629371
$ws.Range("P9").Value = 0.08040263079629371
$ws.Range("Q9").Value = 281.6735961331164
$ws.Range("R9").Value = 2535.062365198048
$ws.Range("S9").Value = 0.02420606264760553
$ws.Range("T9").Value = 0.02420606264760553
$ws.Range("G10").Value = 19.67155566666667
$ws.Range("H10").Value = 59.014667
$ws.Range("I10").Value = 0.2257411748281949
$ws.Range("J10").Value = 0.2257411748281949
$ws.Range("M10").Value = 70.46291600000001
$ws.Range("N10").Value = 211.388748
$ws.Range("O10").Value = 0.5276750397950939
$ws.Range("P10").Value = 0.5276750397950939
$ws.Range("Q10").Value = 1386.115174529657
$ws.Range("R10").Value = 12475.03657076692
$ws.Range("S10").Value = 0.119117983410859
$ws.Range("T10").Value = 0.119117983410859
$ws.Range("G11").Value = 19.67155566666667
$ws.Range("H11").Value = 59.014667
$ws.Range("I11").Value = 0.2257411748281949
$ws.Range("J11").Value = 0.2257411748281949
$ws.Range("O11").Value = 0.07361176802536967
$ws.Range("P11").Value = 0.07361176802536967
$ws.Range("Q11").Value = 193.3659562967841
$ws.Range("R11").Value = 1740.293606671057
$ws.Range("S11").Value = 0.0166172069952275
$ws.Range("T11").Value = 0.0166172069952275
$ws.Range("G12").Value = 19.67155566666667
$ws.Range("H12").Value = 59.014667
$ws.Range("I12").Value = 0.2257411748281949
$ws.Range("J12").Value = 0.2257411748281949
$ws.Range("M12").Value = 42.505498
$ws.Range("N12").Value = 127.516494
$ws.Range("O12").Value = 0.3183105613832428
$ws.Range("P12").Value = 0.3183105613832428
$ws.Range("Q12").Value = 836.1492700463886
$ws.Range("R12").Value = 7525.343430417497
$ws.Range("S12").Value = 0.07185580008687549
$ws.Range("T12").Value = 0.07185580008687549
$ws.Range("G13").Value = 19.67155566666667
$ws.Range("H13").Value = 59.014667
$ws.Range("I13").Value = 0.2257411748281949
$ws.Range("J13").Value = 0.2257411748281949
$ws.Range("M13").Value = 10.73653933333333
$ws.Range("N13").Value = 32.209618
$ws.Range("O13").Value = 0.08040263079629371
$ws.Range("P13").Value = 0.08040263079629371
$ws.Range("Q13").Value = 211.2044311630229
$ws.Range("R13").Value = 1900.839880467206
$ws.Range("S13").Value = 0.01815018433523295
$ws.Range("T13").Value = 0.01815018433523295
$ws.Range("G14").Value = 21.69639766666667
$ws.Range("H14").Value = 65.08919299999999
$ws.Range("I14").Value = 0.2489772736740025
$ws.Range("J14").Value = 0.2489772736740025
$ws.Range("M14").Value = 70.46291600000001
$ws.Range("N14").Value = 211.388748
$ws.Range("O14").Value = 0.5276750397950939
$ws.Range("P14").Value = 0.5276750397950939
$ws.Range("Q14").Value = 1528.791446288929
$ws.Range("R14").Value = 13759.12301660036
$ws.Range("S14").Value = 0.1313790927940032
$ws.Range("T14").Value = 0.1313790927940032
$ws.Range("G15").Value = 21.69639766666667
$ws.Range("H15").Value = 65.08919299999999
$ws.Range("I15").Value = 0.2489772736740025
$ws.Range("J15").Value = 0.2489772736740025
$ws.Range("O15").Value = 0.07361176802536967
$ws.Range("P15").Value = 0.07361176802536967
$ws.Range("Q15").Value = 213.2695936254448
$ws.Range("R15").Value = 1919.426342629003
$ws.Range("S15").Value = 0.01832765731327965
$ws.Range("T15").Value = 0.01832765731327965
$ws.Range("G16").Value = 21.69639766666667
$ws.Range("H16").Value = 65.08919299999999
$ws.Range("I16").Value = 0.2489772736740025
$ws.Range("J16").Value = 0.2489772736740025
$ws.Range("M16").Value = 42.505498
$ws.Range("N16").Value = 127.516494
$ws.Range("O16").Value = 0.3183105613832428
$ws.Range("P16").Value = 0.3183105613832428
$ws.Range("Q16").Value = 922.2161876277046
$ws.Range("R16").Value = 8299.945688649341
$ws.Range("S16").Value = 0.079252095754841
$ws.Range("T16").Value = 0.079252095754841
$ws.Range("G17").Value = 21.69639766666667
$ws.Range("H17").Value = 65.08919299999999
$ws.Range("I17").Value = 0.2489772736740025
$ws.Range("J17").Value = 0.2489772736740025
$ws.Range("M17").Value = 10.73653933333333
$ws.Range("N17").Value = 32.209618
$ws.Range("O17").Value = 0.08040263079629371
$ws.Range("P17").Value = 0.08040263079629371
$ws.Range("Q17").Value = 232.9442269398082
$ws.Range("R17").Value = 2096.498042458274
$ws.Range("S17").Value = 0.0200184278118786
$ws.Range("T17").Value = 0.0200184278118786
